$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix up the "Meaning" column (E) labels that were incorrectly re-using
# other rows' descriptions.
$ws.Range("E10").Value = "MTTR for AC OHL (extreme events)"
$ws.Range("E14").Value = "MTTR for DC OHL (extreme events)"
$ws.Range("E17").Value = "MTTR for DC cable"
$ws.Range("E18").Value = "MTTR for DC cable (extreme events)"

# Move / restore the active selection to E19.
$ws.Range("E19").Select()
